$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new week of price data arrived (date serial 44509). It needs to be
# inserted as the newest 3 records for "Lechuga" (rows 399-401), pushing
# the existing historical rows (399-412) down to (402-415).
$ws.Rows("399:401").Insert()

# --- Row 399: Lechuga / Conconina(o) ---
$ws.Range("A399").Value = 11
$ws.Range("B399").Value = "Vega Monumental Concepción"
$ws.Range("C399").Value = "Bíobío"
$ws.Range("D399").Value = 44509
$ws.Range("E399").Value = 8
$ws.Range("F399").Value = 100112033
$ws.Range("G399").Value = "Lechuga"
$ws.Range("H399").Value = "Conconina(o)"
$ws.Range("I399").Value = "Primera"
$ws.Range("J399").Value = 100
$ws.Range("K399").Value = 5500
$ws.Range("L399").Value = 6000
$ws.Range("M399").Value = 5750
$ws.Range("N399").Value = "$/caja 10 unidades"
$ws.Range("O399").Value = "Región Metropolitana"
$ws.Range("P399").Value = 575
$ws.Range("Q399").Value = 10
$ws.Range("R399").Value = "Hortaliza"

# --- Row 400: Lechuga / Escarola ---
$ws.Range("A400").Value = 11
$ws.Range("B400").Value = "Vega Monumental Concepción"
$ws.Range("C400").Value = "Bíobío"
$ws.Range("D400").Value = 44509
$ws.Range("E400").Value = 8
$ws.Range("F400").Value = 100112033
$ws.Range("G400").Value = "Lechuga"
$ws.Range("H400").Value = "Escarola"
$ws.Range("I400").Value = "Primera"
$ws.Range("J400").Value = 200
$ws.Range("K400").Value = 6500
$ws.Range("L400").Value = 7000
$ws.Range("M400").Value = 6750
$ws.Range("N400").Value = "$/caja 15 unidades"
$ws.Range("O400").Value = "Región de Coquimbo"
$ws.Range("P400").Value = 450
$ws.Range("Q400").Value = 15
$ws.Range("R400").Value = "Hortaliza"

# --- Row 401: Lechuga / Milanesa ---
$ws.Range("A401").Value = 11
$ws.Range("B401").Value = "Vega Monumental Concepción"
$ws.Range("C401").Value = "Bíobío"
$ws.Range("D401").Value = 44509
$ws.Range("E401").Value = 8
$ws.Range("F401").Value = 100112033
$ws.Range("G401").Value = "Lechuga"
$ws.Range("H401").Value = "Milanesa"
$ws.Range("I401").Value = "Primera"
$ws.Range("J401").Value = 100
$ws.Range("K401").Value = 5000
$ws.Range("L401").Value = 5500
$ws.Range("M401").Value = 5250
$ws.Range("N401").Value = "$/caja 20 unidades"
$ws.Range("O401").Value = "Región Metropolitana"
$ws.Range("P401").Value = 262
$ws.Range("Q401").Value = 20
$ws.Range("R401").Value = "Hortaliza"
